$d = $word.ActiveDocument

# The heading "Mô tả chi tiết các đối tượng người dùng" drops the
# word "chi tiết " so it reads "Mô tả các đối tượng người dùng".
# Locate the paragraph's range explicitly (rather than a document-wide
# Find) so we only ever touch this one heading and inherit its existing
# run formatting (Times New Roman / bold / vi-VN) for the replacement
# text automatically.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Mô tả chi tiết các đối tượng người dùng") {
        $target = $p.Range
        break
    }
}

if ($target -ne $null) {
    $found = $target.Find.Execute("Mô tả chi tiết các đối tượng người dùng", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "Mô tả các đối tượng người dùng", 2)
} else {
    # Fallback: scope-free replace (still only matches this unique heading).
    $found = $d.Content.Find.Execute("Mô tả chi tiết các đối tượng người dùng", $true, $false, $false, $false, $false, `
                                      $true, 1, $false, "Mô tả các đối tượng người dùng", 2)
}

Write-Host "Heading updated:" $found
